# Apply "adding new progress as of date 04 nov 2025" update
# For rows 3-18 on the "Training Dashboard" sheet:
#   - Column H (PERIOD TO EXPIRE) decreases by 1
#   - Column I (LAST UPDATE) changes from "03-Nov-2025" text to "04-Nov-2025" text
#
# Note: assigning a date-like string directly to Value/Value2/Formula causes
# Excel to auto-convert it into a real date serial number (and change the
# cell's number format). To keep the cell as literal text (matching the
# original inline-string storage), we build the text via a TEXT() formula,
# then copy/paste-special as values so the literal string is kept without
# triggering date auto-detection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 18; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H - PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # Column I - LAST UPDATE

    # Decrement the numeric "period to expire" value by 1
    $hCell.Value2 = $hCell.Value2 - 1

    # Replace the last-update date text with the new date, keeping it as text
    $iCell.Formula = '=TEXT("04-Nov-2025","@")'
    $iCell.Copy()
    $iCell.PasteSpecial(-4163)  # xlPasteValues
}

$excel.CutCopyMode = 0
